$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: A2 becomes a real number (was a text "08")
$ws.Range("A2").Value = 8

# New data rows 3-6 (Roll Number as plain numbers)
$data = @(
    @(14,    "Traditional OCR Model", 96.73469387755102,  65.95618709295441,  83.79454584299766,  63.7342908438061,   114.4815117120743),
    @(15,    "Traditional OCR Model", 443.4782608695652,  287.1194379391101,  99.71320461609024,  92.54609650843469,  139.7256702184677),
    @(21,    "Traditional OCR Model", 95.28061224489795,  85.27951651197928,  98.53438244339399,  85.14851485148515,  133.6782664060593),
    @(23136, "Traditional OCR Model", 109.5588235294118,  78.42076798269335,  98.05523582977716,  88.27708703374778,  131.9301843643188)
)

$r = 3
foreach ($row in $data) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
    $ws.Cells.Item($r, 7).Value = $row[6]
    $r = $r + 1
}

# Row 7: Roll Number stored as text "23138" (quote-prefixed, like the source data)
$ws.Range("A7").Value = "'23138"
$ws.Range("B7").Value = "Traditional OCR Model"
$ws.Range("C7").Value = 94.84126984126983
$ws.Range("D7").Value = 80.74800735744941
$ws.Range("E7").Value = 98.82107278317123
$ws.Range("F7").Value = 89.0625
$ws.Range("G7").Value = 140.3848469257355
